$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (columns B:E)
$ws.Range("B2").Value = 24.873432544317609
$ws.Range("C2").Value = 16.490162510701225
$ws.Range("D2").Value = 24.797060122536578
$ws.Range("E2").Value = 25.250462001324593

# Row 3 data values (columns B:E)
$ws.Range("B3").Value = 14.383597367489955
$ws.Range("C3").Value = 22.605253653623379
$ws.Range("D3").Value = 17.107705943601673
$ws.Range("E3").Value = 24.66254717183233

# Update selection to reflect the narrower range B1:E3
$ws.Range("B1:E3").Select()
